$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1 with the same text content as described in the diff
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the formatting of the existing header (E1) onto the new headers
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New boolean data cells F2:H4
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false
